$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is Yvette Guerrier (Maire Ferrier). Fill in her email address as a
# clickable mailto hyperlink, the same way Raymond Joanel's email (C6) is
# already set up.
$ws.Range("C4").Value = "guerrierivette@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:guerrierivette@gmail.com")

# Re-apply the exact look of the existing hyperlink cell (border + hyperlink
# font) so the new cell matches the sheet's existing formatting.
$ws.Range("C6").Copy()
$ws.Range("C4").PasteSpecial(-4122)

# C5 (Henry Claude Crepin's email column) also picks up the same bordered
# hyperlink-column look even though no address was filled in for him.
$ws.Range("C6").Copy()
$ws.Range("C5").PasteSpecial(-4122)

# Now that she has an email on file, mark her as having an mWater account,
# matching the "yes" styling already used for row 6.
$ws.Range("E4").Value = "yes"
$ws.Range("E6").Copy()
$ws.Range("E4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection to A3, as left by the editor.
$ws.Range("A3").Select()
